# "se realizo todas las preguntas" - mark attendance ("p" = presente) for the
# remaining dates (columns O/P) for every student except row 9, and widen the
# email column now that more of it needs to stay visible.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New attendance marks ("p", reusing the existing shared string) - one extra
# day-column filled in for every student row except row 9.
$newMarks = @(
    "O3", "P3",
    "O4",
    "P5",
    "P6",
    "P7",
    "P8",
    "P10",
    "P11",
    "O12", "P12",
    "P13",
    "P14"
)

foreach ($addr in $newMarks) {
    $ws.Range($addr).Value = "p"
}

# Column D (email addresses) got widened.
$ws.Columns.Item(4).ColumnWidth = 32.5

# Move/leave the active selection on P9.
$ws.Range("P9").Select()

# Sheet-tab area was resized too (best effort - not all window chrome
# round-trips through the exporter).
$excel.ActiveWindow.TabRatio = 0.614
